{"js": "// Replace the two-digit-division answer strings throughout the document's\n// table cells. Each original \"a\u00f7b=c, d\" string is unique in the document,\n// so a direct search + replace per pair is safe and precise.\nconst replacements = [\n  [\"53\u00f79=5, 8\", \"14\u00f73=4, 2\"],\n  [\"19\u00f76=3, 1\", \"15\u00f77=2, 1\"],\n  [\"24\u00f79=2, 6\", \"39\u00f77=5, 4\"],\n  [\"46\u00f72=23, 0\", \"98\u00f75=19, 3\"],\n  [\"24\u00f75=4, 4\", \"49\u00f79=5, 4\"],\n  [\"52\u00f74=13, 0\", \"14\u00f72=7, 0\"],\n  [\"45\u00f78=5, 5\", \"34\u00f73=11, 1\"],\n  [\"59\u00f72=29, 1\", \"70\u00f74=17, 2\"],\n  [\"55\u00f74=13, 3\", \"67\u00f75=13, 2\"],\n  [\"40\u00f74=10, 0\", \"98\u00f72=49, 0\"],\n  [\"10\u00f72=5, 0\", \"46\u00f75=9, 1\"],\n  [\"27\u00f78=3, 3\", \"31\u00f72=15, 1\"],\n  [\"26\u00f73=8, 2\", \"77\u00f74=19, 1\"],\n  [\"66\u00f79=7, 3\", \"15\u00f79=1, 6\"],\n  [\"96\u00f79=10, 6\", \"17\u00f76=2, 5\"],\n  [\"53\u00f75=10, 3\", \"12\u00f72=6, 0\"],\n  [\"77\u00f77=11, 0\", \"30\u00f79=3, 3\"],\n  [\"97\u00f78=12, 1\", \"28\u00f78=3, 4\"],\n  [\"97\u00f72=48, 1\", \"15\u00f74=3, 3\"],\n  [\"43\u00f79=4, 7\", \"40\u00f75=8, 0\"],\n  [\"91\u00f76=15, 1\", \"55\u00f78=6, 7\"],\n  [\"29\u00f72=14, 1\", \"61\u00f73=20, 1\"],\n  [\"90\u00f76=15, 0\", \"58\u00f79=6, 4\"],\n  [\"81\u00f79=9, 0\", \"15\u00f76=2, 3\"],\n  [\"90\u00f74=22, 2\", \"87\u00f75=17, 2\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the two-digit-division answer strings throughout the document's\n# table cells. Each original \"a\u00f7b=c, d\" string is unique in the document,\n# so a direct Find/Replace per pair is precise and will not clobber other\n# cells.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"53\u00f79=5, 8\", \"14\u00f73=4, 2\"),\n    @(\"19\u00f76=3, 1\", \"15\u00f77=2, 1\"),\n    @(\"24\u00f79=2, 6\", \"39\u00f77=5, 4\"),\n    @(\"46\u00f72=23, 0\", \"98\u00f75=19, 3\"),\n    @(\"24\u00f75=4, 4\", \"49\u00f79=5, 4\"),\n    @(\"52\u00f74=13, 0\", \"14\u00f72=7, 0\"),\n    @(\"45\u00f78=5, 5\", \"34\u00f73=11, 1\"),\n    @(\"59\u00f72=29, 1\", \"70\u00f74=17, 2\"),\n    @(\"55\u00f74=13, 3\", \"67\u00f75=13, 2\"),\n    @(\"40\u00f74=10, 0\", \"98\u00f72=49, 0\"),\n    @(\"10\u00f72=5, 0\", \"46\u00f75=9, 1\"),\n    @(\"27\u00f78=3, 3\", \"31\u00f72=15, 1\"),\n    @(\"26\u00f73=8, 2\", \"77\u00f74=19, 1\"),\n    @(\"66\u00f79=7, 3\", \"15\u00f79=1, 6\"),\n    @(\"96\u00f79=10, 6\", \"17\u00f76=2, 5\"),\n    @(\"53\u00f75=10, 3\", \"12\u00f72=6, 0\"),\n    @(\"77\u00f77=11, 0\", \"30\u00f79=3, 3\"),\n    @(\"97\u00f78=12, 1\", \"28\u00f78=3, 4\"),\n    @(\"97\u00f72=48, 1\", \"15\u00f74=3, 3\"),\n    @(\"43\u00f79=4, 7\", \"40\u00f75=8, 0\"),\n    @(\"91\u00f76=15, 1\", \"55\u00f78=6, 7\"),\n    @(\"29\u00f72=14, 1\", \"61\u00f73=20, 1\"),\n    @(\"90\u00f76=15, 0\", \"58\u00f79=6, 4\"),\n    @(\"81\u00f79=9, 0\", \"15\u00f76=2, 3\"),\n    @(\"90\u00f74=22, 2\", \"87\u00f75=17, 2\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
